# Updated cryptos list values (Price / Volume(1h)) per the source diff.
# D-column values that look numeric are prefixed with a leading
# apostrophe so Excel stores them as literal text (matching the
# original inline-string cell type) instead of coercing to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.732.89"
$ws.Range("E2").Value = "  -2.27%  "
$ws.Range("D3").Value = "3.235.52"
$ws.Range("E3").Value = "  -1.29%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'577.17"
$ws.Range("E5").Value = "  -1.30%  "
$ws.Range("D6").Value = "'172.40"
$ws.Range("E6").Value = "  -3.90%  "
$ws.Range("D7").Value = "'0.627"
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "3.234.14"
$ws.Range("E9").Value = "  -1.32%  "
$ws.Range("E10").Value = "  -2.77%  "
$ws.Range("D11").Value = "'6.77"
$ws.Range("E11").Value = "  +0.38%  "
$ws.Range("D12").Value = "'0.389"
$ws.Range("E12").Value = "  -3.02%  "
$ws.Range("D13").Value = "3.796.69"
$ws.Range("E13").Value = "  -1.33%  "
$ws.Range("E14").Value = "  -3.11%  "
$ws.Range("D15").Value = "64.841.15"
$ws.Range("E15").Value = "  -2.06%  "
$ws.Range("D16").Value = "'25.80"
$ws.Range("E16").Value = "  -2.06%  "
$ws.Range("D17").Value = "3.230.81"
$ws.Range("E17").Value = "  +0.57%  "
$ws.Range("E18").Value = "  -2.95%  "
$ws.Range("D19").Value = "'417.17"
$ws.Range("E19").Value = "  -4.00%  "
$ws.Range("E20").Value = "  -2.35%  "
$ws.Range("D21").Value = "'12.83"
$ws.Range("E21").Value = "  -2.53%  "
$ws.Range("E22").Value = "  -2.33%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").Value = "'70.42"
$ws.Range("E24").Value = "  -1.73%  "
$ws.Range("E25").Value = "  -0.61%  "
$ws.Range("E26").Value = "  +4.23%  "
$ws.Range("E27").Value = "  -1.76%  "
$ws.Range("E28").Value = "  -1.72%  "
$ws.Range("D29").Value = "'9.01"
$ws.Range("E29").Value = "  +1.99%  "
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("E31").Value = "  -4.43%  "
$ws.Range("D32").Value = "'21.84"
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("E34").Value = "  -3.60%  "
$ws.Range("D35").Value = "'6.42"
$ws.Range("E35").Value = "  -2.86%  "
$ws.Range("E36").Value = "  -2.60%  "
$ws.Range("D37").Value = "'157.95"
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("D38").Value = "'1.40"
$ws.Range("E38").Value = "  -1.84%  "
$ws.Range("D39").Value = "2.819.53"
$ws.Range("E39").Value = "  +1.50%  "
$ws.Range("E40").Value = "  -2.91%  "
$ws.Range("D41").Value = "'25.49"
$ws.Range("E41").Value = "  -4.14%  "
$ws.Range("D42").Value = "'4.21"
$ws.Range("E42").Value = "  -2.72%  "
$ws.Range("E43").Value = "  -1.74%  "
$ws.Range("E44").Value = "  -6.34%  "
$ws.Range("D45").Value = "'5.77"
$ws.Range("E45").Value = "  -4.19%  "
$ws.Range("E46").Value = "  -4.25%  "
$ws.Range("D47").Value = "'2.17"
$ws.Range("E47").Value = "  -4.78%  "
$ws.Range("D48").Value = "'302.46"
$ws.Range("E48").Value = "  -5.76%  "
$ws.Range("E49").Value = "  -5.18%  "
$ws.Range("E50").Value = "  -1.32%  "
$ws.Range("E51").Value = "  -1.55%  "
